$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F9").Value = 9
$ws.Range("G9").Value = 3274.92
$ws.Range("F11").Value = 15
$ws.Range("G11").Value = 3268.95
$ws.Range("B12").Value = 7271.82
$ws.Range("F19").Value = 111
$ws.Range("G19").Value = 5686.53
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 138.3
$ws.Range("B27").Value = 11405.64
$ws.Range("F31").Value = 80
$ws.Range("G31").Value = 2498.4
$ws.Range("F36").Value = 39
$ws.Range("G36").Value = 7522.71
$ws.Range("F37").Value = 36
$ws.Range("G37").Value = 1270.44
$ws.Range("F50").Value = 33
$ws.Range("G50").Value = 1165.56
$ws.Range("B56").Value = 46244.99
$ws.Range("F82").Value = 11
$ws.Range("G82").Value = 1553.09
$ws.Range("F91").Value = 118
$ws.Range("G91").Value = 7484.74
$ws.Range("F93").Value = 302
$ws.Range("G93").Value = 19237.4
$ws.Range("F94").Value = 66
$ws.Range("G94").Value = 4641.12
$ws.Range("F104").Value = 153
$ws.Range("G104").Value = 15676.38
$ws.Range("F106").Value = 34
$ws.Range("G106").Value = 4582.18
$ws.Range("F107").Value = 16
$ws.Range("G107").Value = 808.8
$ws.Range("F110").Value = 7
$ws.Range("G110").Value = 1884.33
$ws.Range("F113").Value = 83
$ws.Range("G113").Value = 3925.9
$ws.Range("B115").Value = 261176.55
$ws.Range("F127").Value = 1
$ws.Range("G127").Value = 39.86
$ws.Range("F129").Value = 53
$ws.Range("G129").Value = 5519.95
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0
$ws.Range("B133").Value = 13218.93
$ws.Range("F146").Value = 150
$ws.Range("G146").Value = 2952
$ws.Range("B151").Value = 24854.91
$ws.Range("F158").Value = 9
$ws.Range("G158").Value = 413.1
$ws.Range("B159").Value = 838.15
$ws.Range("F170").Value = 0
$ws.Range("G170").Value = 0
$ws.Range("F173").Value = 3
$ws.Range("G173").Value = 152.22
$ws.Range("B185").Value = 19811.06
$ws.Range("F189").Value = 135
$ws.Range("G189").Value = 4189.05
$ws.Range("B197").Value = 14488.77
$ws.Range("F211").Value = 17
$ws.Range("G211").Value = 817.02
$ws.Range("B212").Value = 951.8
$ws.Range("F222").Value = 42
$ws.Range("G222").Value = 8035.86
$ws.Range("F223").Value = 13
$ws.Range("G223").Value = 1748.11
$ws.Range("B228").Value = 25743.98
$ws.Range("F251").Value = 68
$ws.Range("G251").Value = 4077.96
$ws.Range("B279").Value = 120209.25
$ws.Range("B283").Value = 61610
$ws.Range("D283").Value = 102.71
$ws.Range("E283").Value = 122.71
$ws.Range("F283").Value = 211
$ws.Range("G283").Value = 21671.81
$ws.Range("B284").Value = 57077
$ws.Range("D284").Value = 93.08
$ws.Range("E284").Value = 111.2
$ws.Range("F284").Value = 1
$ws.Range("G284").Value = 93.08
$ws.Range("F319").Value = 0
$ws.Range("G319").Value = 0
$ws.Range("F323").Value = 147
$ws.Range("G323").Value = 14861.7
$ws.Range("F325").Value = 54
$ws.Range("G325").Value = 10723.86
$ws.Range("F330").Value = 117
$ws.Range("G330").Value = 6918.21
$ws.Range("F341").Value = 224
$ws.Range("G341").Value = 15749.44
$ws.Range("B349").Value = 377668.05
$ws.Range("F355").Value = 27
$ws.Range("G355").Value = 9875.790000000001
$ws.Range("F357").Value = 7
$ws.Range("G357").Value = 1527.61
$ws.Range("B358").Value = 31327.08
$ws.Range("F403").Value = 67
$ws.Range("G403").Value = 2493.07
$ws.Range("F404").Value = 67
$ws.Range("G404").Value = 1647.53
$ws.Range("B407").Value = 50937.29
$ws.Range("F423").Value = 71
$ws.Range("G423").Value = 3883.7
$ws.Range("B424").Value = 48890.37
$ws.Range("F435").Value = 626
$ws.Range("G435").Value = 8419.700000000001
$ws.Range("F437").Value = 622
$ws.Range("G437").Value = 7967.82
$ws.Range("F441").Value = 330
$ws.Range("G441").Value = 4227.3
$ws.Range("F442").Value = 388
$ws.Range("G442").Value = 7655.24
$ws.Range("B453").Value = 108472.62
$ws.Range("F463").Value = 1
$ws.Range("G463").Value = 344.57
$ws.Range("F474").Value = 7
$ws.Range("G474").Value = 4087.65
$ws.Range("B478").Value = 42736.07
$ws.Range("F497").Value = 377
$ws.Range("G497").Value = 2582.45
$ws.Range("F500").Value = 430
$ws.Range("G500").Value = 2881
$ws.Range("B504").Value = 41224.06
$ws.Range("F516").Value = 35
$ws.Range("G516").Value = 1826.3
$ws.Range("B525").Value = 28797.4
$ws.Range("F547").Value = 18
$ws.Range("G547").Value = 397.98
$ws.Range("F548").Value = 123
$ws.Range("G548").Value = 5576.82
$ws.Range("B554").Value = 7334.6
$ws.Range("F556").Value = 13
$ws.Range("G556").Value = 5542.94
$ws.Range("F558").Value = 40
$ws.Range("G558").Value = 4154.8
$ws.Range("F563").Value = 55
$ws.Range("G563").Value = 1536.15
$ws.Range("F564").Value = 9
$ws.Range("G564").Value = 246.6
$ws.Range("B571").Value = 44551.17
$ws.Range("F574").Value = 83
$ws.Range("G574").Value = 14775.66
$ws.Range("F577").Value = 99
$ws.Range("G577").Value = 2692.8
$ws.Range("F578").Value = 81
$ws.Range("G578").Value = 2203.2
$ws.Range("F579").Value = 34
$ws.Range("G579").Value = 924.8
$ws.Range("B580").Value = 65767.25999999999
$ws.Range("F599").Value = 51
$ws.Range("G599").Value = 1688.61
$ws.Range("F603").Value = 12
$ws.Range("G603").Value = 397.32
$ws.Range("B608").Value = 33243.52
$ws.Range("F621").Value = 2
$ws.Range("G621").Value = 1361.84
$ws.Range("B627").Value = 13675.15
$ws.Range("F654").Value = 10
$ws.Range("G654").Value = 877
$ws.Range("B655").Value = 3609.4
$ws.Range("F657").Value = 12
$ws.Range("G657").Value = 2828.88
$ws.Range("B664").Value = 22850.78
$ws.Range("F686").Value = 30
$ws.Range("G686").Value = 2357.7
$ws.Range("B694").Value = 24805.92
$ws.Range("F697").Value = 52
$ws.Range("G697").Value = 1944.8
$ws.Range("F701").Value = 100
$ws.Range("G701").Value = 3740
$ws.Range("B702").Value = 10397.74
$ws.Range("F710").Value = 0
$ws.Range("G710").Value = 0
$ws.Range("B714").Value = 31044.68
$ws.Range("F747").Value = 1152
$ws.Range("G747").Value = 187902.72
$ws.Range("F750").Value = 47
$ws.Range("G750").Value = 6949.42
$ws.Range("F751").Value = 116
$ws.Range("G751").Value = 7830
$ws.Range("B752").Value = 220064.4
$ws.Range("B753").Value = 2306453.79
$ws.Range("B754").Value = 2306453.79
